$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row with ticker "GRT-USD" right after the last used row (A62 -> A63)
$ws.Range("A63").Value = "GRT-USD"
